$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add Wins / Losses / Ties columns (AD1:AF1) ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header styling (bold, bordered, centered) used by the other
# header cells by copying the format from the adjacent header cell AC1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows: team record for every player row (2-55) ---
$ws.Range("AD2:AD55").Value = 96
$ws.Range("AE2:AE55").Value = 67
$ws.Range("AF2:AF55").Value = 0
